$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "huevos,harina,vainilla,leche,"
$ws.Range("C3").Value = "huevos,harina,manzana,"
$ws.Range("C4").Value = "huevos,vainilla,harina,"
$ws.Range("C6").Value = "harina, huevos, limon,merengue,crema,"
